$d = $word.ActiveDocument

$d.Content.Find.Execute("87-73=14", $true, $true, $false, $false, $false, $true, 1, $false, "92-83=9", 2) | Out-Null
$d.Content.Find.Execute("39-33=6", $true, $true, $false, $false, $false, $true, 1, $false, "34+48=82", 2) | Out-Null
$d.Content.Find.Execute("32+25=57", $true, $true, $false, $false, $false, $true, 1, $false, "86-76=10", 2) | Out-Null
$d.Content.Find.Execute("76-7=69", $true, $true, $false, $false, $false, $true, 1, $false, "69-12=57", 2) | Out-Null
$d.Content.Find.Execute("32-24=8", $true, $true, $false, $false, $false, $true, 1, $false, "93-63=30", 2) | Out-Null
$d.Content.Find.Execute("78-33=45", $true, $true, $false, $false, $false, $true, 1, $false, "76-39=37", 2) | Out-Null
$d.Content.Find.Execute("45+43=88", $true, $true, $false, $false, $false, $true, 1, $false, "56-31=25", 2) | Out-Null
$d.Content.Find.Execute("50+9=59", $true, $true, $false, $false, $false, $true, 1, $false, "52+17=69", 2) | Out-Null
$d.Content.Find.Execute("33+65=98", $true, $true, $false, $false, $false, $true, 1, $false, "99-85=14", 2) | Out-Null
$d.Content.Find.Execute("70-59=11", $true, $true, $false, $false, $false, $true, 1, $false, "30+2=32", 2) | Out-Null
$d.Content.Find.Execute("49+15=64", $true, $true, $false, $false, $false, $true, 1, $false, "65-33=32", 2) | Out-Null
$d.Content.Find.Execute("11+17=28", $true, $true, $false, $false, $false, $true, 1, $false, "34+60=94", 2) | Out-Null
$d.Content.Find.Execute("36-16=20", $true, $true, $false, $false, $false, $true, 1, $false, "97-93=4", 2) | Out-Null
$d.Content.Find.Execute("0+40=40", $true, $true, $false, $false, $false, $true, 1, $false, "26+45=71", 2) | Out-Null
$d.Content.Find.Execute("14+13=27", $true, $true, $false, $false, $false, $true, 1, $false, "70+5=75", 2) | Out-Null
$d.Content.Find.Execute("34-18=16", $true, $true, $false, $false, $false, $true, 1, $false, "93-16=77", 2) | Out-Null
$d.Content.Find.Execute("29+26=55", $true, $true, $false, $false, $false, $true, 1, $false, "6+13=19", 2) | Out-Null
$d.Content.Find.Execute("5+66=71", $true, $true, $false, $false, $false, $true, 1, $false, "58-57=1", 2) | Out-Null
$d.Content.Find.Execute("71-68=3", $true, $true, $false, $false, $false, $true, 1, $false, "66-8=58", 2) | Out-Null
$d.Content.Find.Execute("19+7=26", $true, $true, $false, $false, $false, $true, 1, $false, "53+2=55", 2) | Out-Null
$d.Content.Find.Execute("48-44=4", $true, $true, $false, $false, $false, $true, 1, $false, "66-47=19", 2) | Out-Null
$d.Content.Find.Execute("36-5=31", $true, $true, $false, $false, $false, $true, 1, $false, "25+49=74", 2) | Out-Null
$d.Content.Find.Execute("38+9=47", $true, $true, $false, $false, $false, $true, 1, $false, "19+35=54", 2) | Out-Null
$d.Content.Find.Execute("56+15=71", $true, $true, $false, $false, $false, $true, 1, $false, "58-2=56", 2) | Out-Null
$d.Content.Find.Execute("23+12=35", $true, $true, $false, $false, $false, $true, 1, $false, "44-10=34", 2) | Out-Null
$d.Content.Find.Execute("69-5=64", $true, $true, $false, $false, $false, $true, 1, $false, "27+16=43", 2) | Out-Null
$d.Content.Find.Execute("35-8=27", $true, $true, $false, $false, $false, $true, 1, $false, "26-7=19", 2) | Out-Null
$d.Content.Find.Execute("39-17=22", $true, $true, $false, $false, $false, $true, 1, $false, "57+22=79", 2) | Out-Null
$d.Content.Find.Execute("75-43=32", $true, $true, $false, $false, $false, $true, 1, $false, "61+12=73", 2) | Out-Null
$d.Content.Find.Execute("86-47=39", $true, $true, $false, $false, $false, $true, 1, $false, "25+55=80", 2) | Out-Null
$d.Content.Find.Execute("88-51=37", $true, $true, $false, $false, $false, $true, 1, $false, "91-72=19", 2) | Out-Null
$d.Content.Find.Execute("38+5=43", $true, $true, $false, $false, $false, $true, 1, $false, "69-56=13", 2) | Out-Null
$d.Content.Find.Execute("32+17=49", $true, $true, $false, $false, $false, $true, 1, $false, "73+18=91", 2) | Out-Null
$d.Content.Find.Execute("93-79=14", $true, $true, $false, $false, $false, $true, 1, $false, "77-23=54", 2) | Out-Null
$d.Content.Find.Execute("83-78=5", $true, $true, $false, $false, $false, $true, 1, $false, "16+82=98", 2) | Out-Null
$d.Content.Find.Execute("95-37=58", $true, $true, $false, $false, $false, $true, 1, $false, "55-28=27", 2) | Out-Null
$d.Content.Find.Execute("67+5=72", $true, $true, $false, $false, $false, $true, 1, $false, "95-63=32", 2) | Out-Null
$d.Content.Find.Execute("38-21=17", $true, $true, $false, $false, $false, $true, 1, $false, "11-3=8", 2) | Out-Null
$d.Content.Find.Execute("22+46=68", $true, $true, $false, $false, $false, $true, 1, $false, "39+54=93", 2) | Out-Null
$d.Content.Find.Execute("61-41=20", $true, $true, $false, $false, $false, $true, 1, $false, "22+11=33", 2) | Out-Null
$d.Content.Find.Execute("82-15=67", $true, $true, $false, $false, $false, $true, 1, $false, "23+11=34", 2) | Out-Null
$d.Content.Find.Execute("14+44=58", $true, $true, $false, $false, $false, $true, 1, $false, "75+16=91", 2) | Out-Null
$d.Content.Find.Execute("69-64=5", $true, $true, $false, $false, $false, $true, 1, $false, "20+70=90", 2) | Out-Null
$d.Content.Find.Execute("44+5=49", $true, $true, $false, $false, $false, $true, 1, $false, "37+33=70", 2) | Out-Null
$d.Content.Find.Execute("46+53=99", $true, $true, $false, $false, $false, $true, 1, $false, "74+9=83", 2) | Out-Null
$d.Content.Find.Execute("45-21=24", $true, $true, $false, $false, $false, $true, 1, $false, "1+88=89", 2) | Out-Null
$d.Content.Find.Execute("55+29=84", $true, $true, $false, $false, $false, $true, 1, $false, "59+23=82", 2) | Out-Null
$d.Content.Find.Execute("63+15=78", $true, $true, $false, $false, $false, $true, 1, $false, "4+56=60", 2) | Out-Null
$d.Content.Find.Execute("25+9=34", $true, $true, $false, $false, $false, $true, 1, $false, "32-12=20", 2) | Out-Null
$d.Content.Find.Execute("86-35=51", $true, $true, $false, $false, $false, $true, 1, $false, "75-71=4", 2) | Out-Null
$d.Content.Find.Execute("66-14=52", $true, $true, $false, $false, $false, $true, 1, $false, "83-76=7", 2) | Out-Null
$d.Content.Find.Execute("79+9=88", $true, $true, $false, $false, $false, $true, 1, $false, "40+42=82", 2) | Out-Null
$d.Content.Find.Execute("89-79=10", $true, $true, $false, $false, $false, $true, 1, $false, "96-4=92", 2) | Out-Null
$d.Content.Find.Execute("31+33=64", $true, $true, $false, $false, $false, $true, 1, $false, "90-71=19", 2) | Out-Null
$d.Content.Find.Execute("32+19=51", $true, $true, $false, $false, $false, $true, 1, $false, "94+2=96", 2) | Out-Null
$d.Content.Find.Execute("79-56=23", $true, $true, $false, $false, $false, $true, 1, $false, "1+17=18", 2) | Out-Null
$d.Content.Find.Execute("79-1=78", $true, $true, $false, $false, $false, $true, 1, $false, "24+23=47", 2) | Out-Null
$d.Content.Find.Execute("77+14=91", $true, $true, $false, $false, $false, $true, 1, $false, "5+7=12", 2) | Out-Null
$d.Content.Find.Execute("83-9=74", $true, $true, $false, $false, $false, $true, 1, $false, "82-20=62", 2) | Out-Null
$d.Content.Find.Execute("42+47=89", $true, $true, $false, $false, $false, $true, 1, $false, "60-54=6", 2) | Out-Null
$d.Content.Find.Execute("37+5=42", $true, $true, $false, $false, $false, $true, 1, $false, "23-16=7", 2) | Out-Null
$d.Content.Find.Execute("6+66=72", $true, $true, $false, $false, $false, $true, 1, $false, "64-38=26", 2) | Out-Null
$d.Content.Find.Execute("17-1=16", $true, $true, $false, $false, $false, $true, 1, $false, "9+15=24", 2) | Out-Null
$d.Content.Find.Execute("3+84=87", $true, $true, $false, $false, $false, $true, 1, $false, "92-75=17", 2) | Out-Null
$d.Content.Find.Execute("35+16=51", $true, $true, $false, $false, $false, $true, 1, $false, "79+19=98", 2) | Out-Null
$d.Content.Find.Execute("85-56=29", $true, $true, $false, $false, $false, $true, 1, $false, "67-56=11", 2) | Out-Null
$d.Content.Find.Execute("74-59=15", $true, $true, $false, $false, $false, $true, 1, $false, "2+66=68", 2) | Out-Null
$d.Content.Find.Execute("10+88=98", $true, $true, $false, $false, $false, $true, 1, $false, "62-3=59", 2) | Out-Null
$d.Content.Find.Execute("36-31=5", $true, $true, $false, $false, $false, $true, 1, $false, "30+49=79", 2) | Out-Null
$d.Content.Find.Execute("88-68=20", $true, $true, $false, $false, $false, $true, 1, $false, "75-61=14", 2) | Out-Null
$d.Content.Find.Execute("34-27=7", $true, $true, $false, $false, $false, $true, 1, $false, "27+31=58", 2) | Out-Null
$d.Content.Find.Execute("80-34=46", $true, $true, $false, $false, $false, $true, 1, $false, "56-49=7", 2) | Out-Null
$d.Content.Find.Execute("45-19=26", $true, $true, $false, $false, $false, $true, 1, $false, "48+6=54", 2) | Out-Null
$d.Content.Find.Execute("51+43=94", $true, $true, $false, $false, $false, $true, 1, $false, "46-7=39", 2) | Out-Null
$d.Content.Find.Execute("47-47=0", $true, $true, $false, $false, $false, $true, 1, $false, "52+14=66", 2) | Out-Null
$d.Content.Find.Execute("98-86=12", $true, $true, $false, $false, $false, $true, 1, $false, "64-57=7", 2) | Out-Null
$d.Content.Find.Execute("81-50=31", $true, $true, $false, $false, $false, $true, 1, $false, "28+32=60", 2) | Out-Null
$d.Content.Find.Execute("83-21=62", $true, $true, $false, $false, $false, $true, 1, $false, "70+19=89", 2) | Out-Null
$d.Content.Find.Execute("87-25=62", $true, $true, $false, $false, $false, $true, 1, $false, "4+93=97", 2) | Out-Null
$d.Content.Find.Execute("82-68=14", $true, $true, $false, $false, $false, $true, 1, $false, "79-2=77", 2) | Out-Null
$d.Content.Find.Execute("26+1=27", $true, $true, $false, $false, $false, $true, 1, $false, "32+67=99", 2) | Out-Null
$d.Content.Find.Execute("40-27=13", $true, $true, $false, $false, $false, $true, 1, $false, "8+58=66", 2) | Out-Null
$d.Content.Find.Execute("92-28=64", $true, $true, $false, $false, $false, $true, 1, $false, "1+87=88", 2) | Out-Null
$d.Content.Find.Execute("38-35=3", $true, $true, $false, $false, $false, $true, 1, $false, "58-30=28", 2) | Out-Null
$d.Content.Find.Execute("30-5=25", $true, $true, $false, $false, $false, $true, 1, $false, "76-5=71", 2) | Out-Null
$d.Content.Find.Execute("2+94=96", $true, $true, $false, $false, $false, $true, 1, $false, "11+82=93", 2) | Out-Null
$d.Content.Find.Execute("32+12=44", $true, $true, $false, $false, $false, $true, 1, $false, "78-52=26", 2) | Out-Null
$d.Content.Find.Execute("66-19=47", $true, $true, $false, $false, $false, $true, 1, $false, "80-73=7", 2) | Out-Null
$d.Content.Find.Execute("85-18=67", $true, $true, $false, $false, $false, $true, 1, $false, "46+49=95", 2) | Out-Null
$d.Content.Find.Execute("13+51=64", $true, $true, $false, $false, $false, $true, 1, $false, "52+13=65", 2) | Out-Null
$d.Content.Find.Execute("65+32=97", $true, $true, $false, $false, $false, $true, 1, $false, "49-44=5", 2) | Out-Null
$d.Content.Find.Execute("41+52=93", $true, $true, $false, $false, $false, $true, 1, $false, "33-6=27", 2) | Out-Null
$d.Content.Find.Execute("63+12=75", $true, $true, $false, $false, $false, $true, 1, $false, "20+15=35", 2) | Out-Null
$d.Content.Find.Execute("56-44=12", $true, $true, $false, $false, $false, $true, 1, $false, "49+38=87", 2) | Out-Null
$d.Content.Find.Execute("46+20=66", $true, $true, $false, $false, $false, $true, 1, $false, "58-29=29", 2) | Out-Null
$d.Content.Find.Execute("47-19=28", $true, $true, $false, $false, $false, $true, 1, $false, "19-13=6", 2) | Out-Null
$d.Content.Find.Execute("57+24=81", $true, $true, $false, $false, $false, $true, 1, $false, "66+3=69", 2) | Out-Null
$d.Content.Find.Execute("81-40=41", $true, $true, $false, $false, $false, $true, 1, $false, "35+14=49", 2) | Out-Null
$d.Content.Find.Execute("21+25=46", $true, $true, $false, $false, $false, $true, 1, $false, "41+40=81", 2) | Out-Null
$d.Content.Find.Execute("9+72=81", $true, $true, $false, $false, $false, $true, 1, $false, "97-41=56", 2) | Out-Null
